$wb = $excel.ActiveWorkbook

# The "meta" sheet holds key/value pairs describing the chart config.
# Insert a new "style" / "default" pair right after the existing
# "line_symbol" row (row 6), pushing the old trailing blank row
# (formerly row 7) down to row 8.
$meta = $wb.Worksheets.Item("meta")

# Move the existing blank, styled row (old A7) down to row 8 first,
# preserving its formatting (copy carries the cell style along).
$meta.Range("A7").Copy($meta.Range("A8"))

# Now populate the freed-up row 7 with the new key/value pair; A7 keeps
# the same bold/orange "key" style it already had.
$meta.Range("A7").Value = "style"
$meta.Range("B7").Value = "default"
